$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44333
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 58

# Row 3
$ws.Range("D3").Value = 44333
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("S3").Value = 900

# Row 4
$ws.Range("D4").Value = 44333
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("S4").Value = 800

# Row 5
$ws.Range("D5").Value = 45096
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1300

# Row 6
$ws.Range("D6").Value = 45096
$ws.Range("M6").Value = 68
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 1200

# Row 7
$ws.Range("D7").Value = 45096
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 60
$ws.Range("R7").Value = 'Región de O''Higgins'

# Row 8
$ws.Range("D8").Value = 44307
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44323
$ws.Range("L9").Value = 'Primera'
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44323
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 9000
$ws.Range("R10").Value = 'Provincia de Quillota'
$ws.Range("S10").Value = 900

# Row 11
$ws.Range("D11").Value = 45091
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 54
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("S11").Value = 1400

# Row 12
$ws.Range("D12").Value = 45091
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 58
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1200

# Row 13
$ws.Range("D13").Value = 45091
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 48
$ws.Range("R13").Value = 'Región de O''Higgins'

# Row 14
$ws.Range("D14").Value = 44302
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 45
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44308
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("R15").Value = 'Provincia de Quillota'
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44308
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 48
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("R16").Value = 'Provincia de Quillota'
$ws.Range("S16").Value = 800

# Row 17
$ws.Range("D17").Value = 44314
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 47
$ws.Range("N17").Value = 9000
$ws.Range("O17").Value = 9000
$ws.Range("P17").Value = 9000
$ws.Range("R17").Value = 'Provincia de Quillota'
$ws.Range("S17").Value = 900

# Row 18
$ws.Range("D18").Value = 44312
$ws.Range("M18").Value = 48

# Row 19
$ws.Range("D19").Value = 44306

# Row 20
$ws.Range("D20").Value = 44326
$ws.Range("M20").Value = 65
$ws.Range("N20").Value = 10000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 10000
$ws.Range("S20").Value = 1000

# Row 21
$ws.Range("D21").Value = 44326
$ws.Range("M21").Value = 67
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 800

# Row 22
$ws.Range("D22").Value = 45099
$ws.Range("L22").Value = 'Especial'
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1200

# Row 23
$ws.Range("D23").Value = 45099
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 68
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 1000

# Row 24
$ws.Range("D24").Value = 45099
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 9000
$ws.Range("P24").Value = 9000
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 900

# Row 25
$ws.Range("D25").Value = 45106
$ws.Range("M25").Value = 56
$ws.Range("N25").Value = 12000
$ws.Range("O25").Value = 12000
$ws.Range("P25").Value = 12000
$ws.Range("S25").Value = 1200

# Row 26
$ws.Range("D26").Value = 45106
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("S26").Value = 1000

# Row 27
$ws.Range("D27").Value = 45106
$ws.Range("M27").Value = 54
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 8000
$ws.Range("S27").Value = 800

# Row 28
$ws.Range("D28").Value = 44315
$ws.Range("M28").Value = 45
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("S28").Value = 1000

# Row 29
$ws.Range("D29").Value = 44322
$ws.Range("M29").Value = 56

# Row 30
$ws.Range("D30").Value = 44322
$ws.Range("L30").Value = 'Segunda'
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = 8000
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 8000
$ws.Range("S30").Value = 800

# Row 31
$ws.Range("D31").Value = 45092
$ws.Range("L31").Value = 'Especial'
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 13000
$ws.Range("O31").Value = 13000
$ws.Range("P31").Value = 13000
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1300

# Row 32
$ws.Range("D32").Value = 45092
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 65
$ws.Range("N32").Value = 12000
$ws.Range("O32").Value = 12000
$ws.Range("P32").Value = 12000
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 1200

# Row 33
$ws.Range("D33").Value = 45092
$ws.Range("L33").Value = 'Segunda'
$ws.Range("N33").Value = 10000
$ws.Range("O33").Value = 10000
$ws.Range("P33").Value = 10000
$ws.Range("S33").Value = 1000

# Row 34
$ws.Range("D34").Value = 45082
$ws.Range("L34").Value = 'Especial'
$ws.Range("M34").Value = 56
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 15000
$ws.Range("S34").Value = 1500

# Row 35
$ws.Range("D35").Value = 45082
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 67
$ws.Range("N35").Value = 12000
$ws.Range("O35").Value = 12000
$ws.Range("P35").Value = 12000
$ws.Range("S35").Value = 1200

# Row 36
$ws.Range("D36").Value = 45082
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 60
$ws.Range("N36").Value = 10000
$ws.Range("O36").Value = 10000
$ws.Range("P36").Value = 10000
$ws.Range("R36").Value = 'Región de O''Higgins'
$ws.Range("S36").Value = 1000

# Row 37
$ws.Range("D37").Value = 44301
$ws.Range("M37").Value = 45

# Row 38
$ws.Range("D38").Value = 44329
$ws.Range("M38").Value = 56
$ws.Range("N38").Value = 9000
$ws.Range("O38").Value = 9000
$ws.Range("P38").Value = 9000
$ws.Range("R38").Value = 'Región Metropolitana'
$ws.Range("S38").Value = 900

# Row 39
$ws.Range("D39").Value = 44329
$ws.Range("M39").Value = 50
$ws.Range("R39").Value = 'Región Metropolitana'

# Row 40
$ws.Range("D40").Value = 44321
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 58
$ws.Range("N40").Value = 9000
$ws.Range("O40").Value = 9000
$ws.Range("P40").Value = 9000
$ws.Range("R40").Value = 'Provincia de Quillota'
$ws.Range("S40").Value = 900

# Row 41
$ws.Range("D41").Value = 44328
$ws.Range("M41").Value = 45
$ws.Range("N41").Value = 8000
$ws.Range("O41").Value = 8000
$ws.Range("P41").Value = 8000
$ws.Range("R41").Value = 'Provincia de Quillota'
$ws.Range("S41").Value = 800

# Row 42
$ws.Range("D42").Value = 44328
$ws.Range("N42").Value = 7000
$ws.Range("O42").Value = 7000
$ws.Range("P42").Value = 7000
$ws.Range("R42").Value = 'Provincia de Quillota'
$ws.Range("S42").Value = 700

# Row 43
$ws.Range("D43").Value = 44309

# Row 44
$ws.Range("D44").Value = 44699
$ws.Range("L44").Value = 'Especial'
$ws.Range("M44").Value = 56
$ws.Range("N44").Value = 12000
$ws.Range("O44").Value = 12000
$ws.Range("P44").Value = 12000
$ws.Range("S44").Value = 1200

# Row 45
$ws.Range("D45").Value = 44699
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 60
$ws.Range("N45").Value = 10000
$ws.Range("O45").Value = 10000
$ws.Range("P45").Value = 10000
$ws.Range("R45").Value = 'Provincia de Quillota'
$ws.Range("S45").Value = 1000

# Row 46
$ws.Range("D46").Value = 44319
$ws.Range("M46").Value = 68
$ws.Range("R46").Value = 'Provincia de Quillota'

# Row 47
$ws.Range("D47").Value = 44319
$ws.Range("M47").Value = 57
$ws.Range("R47").Value = 'Provincia de Quillota'
